$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Run Mode" column (C) for these rows was switched from "Yes" to "No"
$rows = @(2,3,4,5,6,7,8,9,10,20,21,22,23,24,25,26,27,28,29,30,31,32,33)

foreach ($r in $rows) {
    $ws.Range("C$r").Value = "No"
}

# Reflect the resulting selection/scroll position shown in the edited workbook
$ws.Range("C11:C19").Select()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
